# Refresh the cryptocurrency price list (columns D = Price, E = Volume(1h))
# on Sheet1 to the latest scraped snapshot. Only the cells whose value
# actually changed between scrapes are touched; every other cell (coin
# name, link, rank index, unchanged prices/percentages) is left as-is.
#
# A handful of "Price" values look like plain numbers once the old value
# is overwritten (e.g. "503.21"), and Excel's COM layer auto-converts a
# numeric-looking string assigned via .Value into a real number. The
# source data stores these as text (note values like "56.482.16" with two
# dots, which aren't valid numbers anyway), so for just those cells we set
# NumberFormat to Text ("@") first to force the write to stay text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "56.482.16"; DForceText = $false; E = "  -1.54%  " },
    @{ Row = 3; D = "2.380.38"; DForceText = $false; E = "  -1.17%  " },
    @{ Row = 4; D = $null; DForceText = $false; E = "  -0.06%  " },
    @{ Row = 5; D = "503.21"; DForceText = $true; E = "  -0.67%  " },
    @{ Row = 6; D = "130.67"; DForceText = $true; E = "  -1.79%  " },
    @{ Row = 7; D = "0.998"; DForceText = $true; E = "  +0.35%  " },
    @{ Row = 8; D = "0.547"; DForceText = $true; E = "  -1.84%  " },
    @{ Row = 9; D = "2.389.21"; DForceText = $false; E = "  -2.21%  " },
    @{ Row = 10; D = $null; DForceText = $false; E = "  +1.06%  " },
    @{ Row = 11; D = $null; DForceText = $false; E = "  +0.43%  " },
    @{ Row = 12; D = $null; DForceText = $false; E = "  +1.37%  " },
    @{ Row = 13; D = "4.76"; DForceText = $true; E = "  +3.56%  " },
    @{ Row = 14; D = "2.804.88"; DForceText = $false; E = "  -1.51%  " },
    @{ Row = 15; D = "56.443.73"; DForceText = $false; E = "  -1.32%  " },
    @{ Row = 16; D = "21.64"; DForceText = $true; E = "  -1.20%  " },
    @{ Row = 17; D = $null; DForceText = $false; E = "  -0.17%  " },
    @{ Row = 18; D = "2.394.99"; DForceText = $false; E = "  -1.77%  " },
    @{ Row = 19; D = $null; DForceText = $false; E = "  -2.13%  " },
    @{ Row = 20; D = "4.04"; DForceText = $true; E = "  -1.60%  " },
    @{ Row = 21; D = "308.01"; DForceText = $true; E = "  -1.87%  " },
    @{ Row = 22; D = $null; DForceText = $false; E = "  -1.72%  " },
    @{ Row = 23; D = $null; DForceText = $false; E = "  +0.29%  " },
    @{ Row = 24; D = "65.55"; DForceText = $true; E = "  +0.49%  " },
    @{ Row = 25; D = $null; DForceText = $false; E = "  +0.32%  " },
    @{ Row = 26; D = $null; DForceText = $false; E = "  -3.43%  " },
    @{ Row = 27; D = $null; DForceText = $false; E = "  -2.90%  " },
    @{ Row = 28; D = "7.30"; DForceText = $true; E = "  -3.54%  " },
    @{ Row = 29; D = "172.02"; DForceText = $true; E = "  -1.08%  " },
    @{ Row = 30; D = $null; DForceText = $false; E = "  -1.42%  " },
    @{ Row = 31; D = $null; DForceText = $false; E = "  -2.37%  " },
    @{ Row = 32; D = $null; DForceText = $false; E = "  +0.21%  " },
    @{ Row = 33; D = $null; DForceText = $false; E = "  -2.98%  " },
    @{ Row = 34; D = "5.78"; DForceText = $true; E = "  -5.88%  " },
    @{ Row = 35; D = "0.997"; DForceText = $true; E = "  +0.54%  " },
    @{ Row = 36; D = $null; DForceText = $false; E = "  -1.81%  " },
    @{ Row = 37; D = $null; DForceText = $false; E = "  -4.80%  " },
    @{ Row = 38; D = $null; DForceText = $false; E = "  -1.03%  " },
    @{ Row = 39; D = $null; DForceText = $false; E = "  -1.03%  " },
    @{ Row = 40; D = "0.797"; DForceText = $true; E = "  -1.74%  " },
    @{ Row = 41; D = $null; DForceText = $false; E = "  -4.30%  " },
    @{ Row = 42; D = "131.01"; DForceText = $true; E = "  -2.50%  " },
    @{ Row = 43; D = $null; DForceText = $false; E = "  -0.31%  " },
    @{ Row = 44; D = $null; DForceText = $false; E = "  -0.55%  " },
    @{ Row = 45; D = "0.566"; DForceText = $true; E = "  -0.75%  " },
    @{ Row = 46; D = $null; DForceText = $false; E = "  -0.68%  " },
    @{ Row = 47; D = "242.34"; DForceText = $true; E = "  -5.52%  " },
    @{ Row = 48; D = "0.0484"; DForceText = $true; E = "  -1.76%  " },
    @{ Row = 49; D = $null; DForceText = $false; E = "  -2.02%  " },
    @{ Row = 50; D = "17.19"; DForceText = $true; E = "  +0.74%  " },
    @{ Row = 51; D = $null; DForceText = $false; E = "  -2.02%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cellD = $ws.Range("D$row")
        if ($u.DForceText) {
            $cellD.NumberFormat = "@"
        }
        $cellD.Value = $u.D
    }

    $ws.Range("E$row").Value = $u.E
}
